# Applies the "discount" update to the invoice workbook:
#   - Client discount amount (E18) changes from 1 to 125 (the AMOUNT column
#     F18 recalculates automatically via its shared formula, as do the
#     downstream SUBTOTAL/TAX/TOTAL cells F21, F23, F24).
#   - The footer contact placeholder text (A31) is updated from
#     "name, email address" to the actual client email address.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the client discount amount.
$ws.Range("E18").Value = 125

# Update the footer note with the real contact email address.
$ws.Range("A31").Value = "charlie.charlie@mail.com"
